# Scheduled market-data refresh: update Leve profit sheets with latest
# currentAveragePrice / LevePrice / LeveProfit figures per sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 354.68
$ws.Cells.Item(15, 9).Value = 354.68
$ws.Cells.Item(15, 11).Value = 1064.04
$ws.Cells.Item(15, 13).Value = -895.04
# Row 53
$ws.Cells.Item(53, 8).Value = 295.52942
$ws.Cells.Item(53, 9).Value = 233.125
$ws.Cells.Item(53, 10).Value = 351
$ws.Cells.Item(53, 11).Value = 233.125
$ws.Cells.Item(53, 12).Value = 351
$ws.Cells.Item(53, 13).Value = 403.875
$ws.Cells.Item(53, 14).Value = -1625
# Row 62
$ws.Cells.Item(62, 8).Value = 10406.546
$ws.Cells.Item(62, 9).Value = 11037.2
$ws.Cells.Item(62, 10).Value = 4100
$ws.Cells.Item(62, 11).Value = 11037.2
$ws.Cells.Item(62, 12).Value = 4100
$ws.Cells.Item(62, 13).Value = -10413.2
$ws.Cells.Item(62, 14).Value = -5348
# Row 64
$ws.Cells.Item(64, 8).Value = 18607.303
$ws.Cells.Item(64, 9).Value = 3330.7896
$ws.Cells.Item(64, 10).Value = 39339.715
$ws.Cells.Item(64, 11).Value = 3330.7896
$ws.Cells.Item(64, 12).Value = 39339.715
$ws.Cells.Item(64, 13).Value = -3082.7896
$ws.Cells.Item(64, 14).Value = -39835.715
# Row 65
$ws.Cells.Item(65, 8).Value = 10406.546
$ws.Cells.Item(65, 9).Value = 11037.2
$ws.Cells.Item(65, 10).Value = 4100
$ws.Cells.Item(65, 11).Value = 55186
$ws.Cells.Item(65, 12).Value = 20500
$ws.Cells.Item(65, 13).Value = -52066
$ws.Cells.Item(65, 14).Value = -26740
# Row 67
$ws.Cells.Item(67, 8).Value = 18607.303
$ws.Cells.Item(67, 9).Value = 3330.7896
$ws.Cells.Item(67, 10).Value = 39339.715
$ws.Cells.Item(67, 11).Value = 3330.7896
$ws.Cells.Item(67, 12).Value = 39339.715
$ws.Cells.Item(67, 13).Value = -2472.7896
$ws.Cells.Item(67, 14).Value = -41055.715
# Row 103
$ws.Cells.Item(103, 8).Value = 825
$ws.Cells.Item(103, 9).Value = 697.2222
$ws.Cells.Item(103, 10).Value = 907.1429000000001
$ws.Cells.Item(103, 11).Value = 2091.6666
$ws.Cells.Item(103, 12).Value = 2721.4287
$ws.Cells.Item(103, 13).Value = -1505.6666
$ws.Cells.Item(103, 14).Value = -3893.4287
# Row 132
$ws.Cells.Item(132, 8).Value = 205745.16
$ws.Cells.Item(132, 9).Value = 1573.9762
$ws.Cells.Item(132, 10).Value = 1430772.2
$ws.Cells.Item(132, 11).Value = 4721.9286
$ws.Cells.Item(132, 12).Value = 4292316.6
$ws.Cells.Item(132, 13).Value = -2191.9286
$ws.Cells.Item(132, 14).Value = -4297376.6
# Row 134
$ws.Cells.Item(134, 8).Value = 39396.668
$ws.Cells.Item(134, 10).Value = 39396.668
$ws.Cells.Item(134, 12).Value = 39396.668
$ws.Cells.Item(134, 14).Value = -49536.668
# Row 135
$ws.Cells.Item(135, 8).Value = 20834242
$ws.Cells.Item(135, 9).Value = 605.7222
$ws.Cells.Item(135, 10).Value = 83335150
$ws.Cells.Item(135, 11).Value = 5451.499800000001
$ws.Cells.Item(135, 12).Value = 750016350
$ws.Cells.Item(135, 13).Value = -2916.499800000001
$ws.Cells.Item(135, 14).Value = -750021420
# Row 137
$ws.Cells.Item(137, 8).Value = 4923.6665
$ws.Cells.Item(137, 9).Value = 936.1429000000001
$ws.Cells.Item(137, 10).Value = 6319.3
$ws.Cells.Item(137, 11).Value = 2808.4287
$ws.Cells.Item(137, 12).Value = 18957.9
$ws.Cells.Item(137, 13).Value = -258.4287000000004
$ws.Cells.Item(137, 14).Value = -24057.9
# Row 138
$ws.Cells.Item(138, 8).Value = 1643.62
$ws.Cells.Item(138, 9).Value = 963.8
$ws.Cells.Item(138, 10).Value = 2096.8333
$ws.Cells.Item(138, 11).Value = 2891.4
$ws.Cells.Item(138, 12).Value = 6290.499899999999
$ws.Cells.Item(138, 13).Value = 2248.6
$ws.Cells.Item(138, 14).Value = -16570.4999
# Row 141
$ws.Cells.Item(141, 8).Value = 2337.7144
$ws.Cells.Item(141, 9).Value = 1879.2
$ws.Cells.Item(141, 10).Value = 3484
$ws.Cells.Item(141, 11).Value = 5637.6
$ws.Cells.Item(141, 12).Value = 10452
$ws.Cells.Item(141, 13).Value = -457.6000000000004
$ws.Cells.Item(141, 14).Value = -20812

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Cells.Item(2, 8).Value = 5497.227
$ws.Cells.Item(2, 9).Value = 1051.8334
$ws.Cells.Item(2, 10).Value = 25501.5
$ws.Cells.Item(2, 11).Value = 1051.8334
$ws.Cells.Item(2, 12).Value = 25501.5
$ws.Cells.Item(2, 13).Value = -938.8334
$ws.Cells.Item(2, 14).Value = -25727.5
# Row 116
$ws.Cells.Item(116, 8).Value = 5497.227
$ws.Cells.Item(116, 9).Value = 1051.8334
$ws.Cells.Item(116, 10).Value = 25501.5
$ws.Cells.Item(116, 11).Value = 1051.8334
$ws.Cells.Item(116, 12).Value = 25501.5
$ws.Cells.Item(116, 13).Value = 1242.1666
$ws.Cells.Item(116, 14).Value = -30089.5
# Row 122
$ws.Cells.Item(122, 8).Value = 1176.125
$ws.Cells.Item(122, 9).Value = 1150.3334
$ws.Cells.Item(122, 10).Value = 1253.5
$ws.Cells.Item(122, 11).Value = 3451.0002
$ws.Cells.Item(122, 12).Value = 3760.5
$ws.Cells.Item(122, 13).Value = -1001.0002
$ws.Cells.Item(122, 14).Value = -8660.5

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Cells.Item(3, 8).Value = 5497.227
$ws.Cells.Item(3, 9).Value = 1051.8334
$ws.Cells.Item(3, 10).Value = 25501.5
$ws.Cells.Item(3, 11).Value = 1051.8334
$ws.Cells.Item(3, 12).Value = 25501.5
$ws.Cells.Item(3, 13).Value = -937.8334
$ws.Cells.Item(3, 14).Value = -25729.5
# Row 134
$ws.Cells.Item(134, 8).Value = 43524220
$ws.Cells.Item(134, 9).Value = 1691.6
$ws.Cells.Item(134, 10).Value = 77003090
$ws.Cells.Item(134, 11).Value = 5074.799999999999
$ws.Cells.Item(134, 12).Value = 231009270
$ws.Cells.Item(134, 13).Value = -2539.799999999999
$ws.Cells.Item(134, 14).Value = -231014340

$ws = $wb.Worksheets.Item("CRP")
# Row 116
$ws.Cells.Item(116, 8).Value = 22362.857
$ws.Cells.Item(116, 10).Value = 22362.857
$ws.Cells.Item(116, 12).Value = 22362.857
$ws.Cells.Item(116, 14).Value = -31540.857

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Cells.Item(113, 8).Value = 415.64706
$ws.Cells.Item(113, 9).Value = 400.375
$ws.Cells.Item(113, 10).Value = 429.22223
$ws.Cells.Item(113, 11).Value = 1201.125
$ws.Cells.Item(113, 12).Value = 1287.66669
$ws.Cells.Item(113, 13).Value = 968.875
$ws.Cells.Item(113, 14).Value = -5627.66669

$ws = $wb.Worksheets.Item("GSM")
# Row 140
$ws.Cells.Item(140, 8).Value = 57419.2
$ws.Cells.Item(140, 10).Value = 57419.2
$ws.Cells.Item(140, 12).Value = 57419.2
$ws.Cells.Item(140, 14).Value = -67779.2

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 2056.2727
$ws.Cells.Item(68, 9).Value = 1366.5
$ws.Cells.Item(68, 10).Value = 2884
$ws.Cells.Item(68, 11).Value = 1366.5
$ws.Cells.Item(68, 12).Value = 2884
$ws.Cells.Item(68, 13).Value = -617.5
$ws.Cells.Item(68, 14).Value = -4382
# Row 71
$ws.Cells.Item(71, 8).Value = 2056.2727
$ws.Cells.Item(71, 9).Value = 1366.5
$ws.Cells.Item(71, 10).Value = 2884
$ws.Cells.Item(71, 11).Value = 6832.5
$ws.Cells.Item(71, 12).Value = 14420
$ws.Cells.Item(71, 13).Value = -3088.5
$ws.Cells.Item(71, 14).Value = -21908
# Row 122
$ws.Cells.Item(122, 8).Value = 2257.3447
$ws.Cells.Item(122, 9).Value = 2200.65
$ws.Cells.Item(122, 10).Value = 2383.3333
$ws.Cells.Item(122, 11).Value = 6601.950000000001
$ws.Cells.Item(122, 12).Value = 7149.999899999999
$ws.Cells.Item(122, 13).Value = -4151.950000000001
$ws.Cells.Item(122, 14).Value = -12049.9999

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 3612.439
$ws.Cells.Item(132, 9).Value = 784.06665
$ws.Cells.Item(132, 10).Value = 11326.182
$ws.Cells.Item(132, 11).Value = 2352.19995
$ws.Cells.Item(132, 12).Value = 33978.546
$ws.Cells.Item(132, 13).Value = 177.8000499999998
$ws.Cells.Item(132, 14).Value = -39038.546
